# Remove the "Footer component" bullet paragraph entirely.
$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`n", [char]7) -eq "Footer component") {
        $p.Range.Delete()
        break
    }
}
